# Apply crypto price / coin-rank refresh per commit "Updated symbol list on Mon Dec 26 09:20:12 UTC 2022"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds numeric-looking values that the source workbook stores as TEXT
# (inline strings, e.g. "243.02"). A plain Range.Value assignment would let Excel
# auto-convert these into real numbers, which would change the cell type. To keep them
# as text (matching the source data) without leaving a stray number-format override
# behind, each cell is: (1) switched to Text format, (2) given its new text value, then
# (3) reset back to the default "Normal" style so formatting is unchanged overall.
$priceUpdates = @{
    "D2" = '242.93'
    "D3" = '23.00'
    "D4" = '5.411'
    "D5" = '0.05922'
    "D6" = '3.451'
    "D7" = '6.541'
    "D8" = '0.8104'
    "D9" = '0.9098'
    "D10" = '0.0005941'
    "D11" = '0.1402'
    "D12" = '0.07373'
    "D13" = '0.03265'
    "D14" = '0.03044'
    "D15" = '0.09348'
    "D16" = '3.852'
    "D17" = '0.001582'
    "D18" = '0.04671'
    "D19" = '0.006049'
    "D20" = '0.004972'
    "D21" = '0.0009870'
    "D22" = '0.00009102'
    "D24" = '2.138'
    "D25" = '0.3239'
    "D40" = '0.03960'
    "D42" = '0.1075'
    "D43" = '0.003001'
    "D44" = '0.008015'
    "D45" = '0.00005246'
    "D47" = '0.7822'
    "D48" = '0.002270'
    "D49" = '0.00002100'
    "D50" = '0.0002000'
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}

# Coin name / link / volume-label columns are plain (non-numeric-looking) text,
# so they can be written directly without the text-format dance above.
$textUpdates = @{
    "B10" = 'One'
    "C10" = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    "E10" = '9OneONE'
    "B11" = 'WazirX'
    "C11" = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    "E11" = '10WazirXWRX'
    "B12" = 'MandalaExchangeToken'
    "C12" = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    "E12" = '11MandalaExchangeTokenMDX'
    "B13" = 'LiechtensteinCryptoassetsExchange'
    "C13" = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    "E13" = '12LiechtensteinCryptoassetsExchangeLCX'
    "B14" = 'BitrueCoin'
    "C14" = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    "E14" = '13BitrueCoinBTR'
    "B15" = 'BitMartToken'
    "C15" = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    "E15" = '14BitMartTokenBMX'
    "B16" = 'MCDex'
    "C16" = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    "E16" = '15MCDexMCB'
    "B17" = 'BitForexToken'
    "C17" = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    "E17" = '16BitForexTokenBF'
    "B18" = 'CoinExToken'
    "C18" = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    "E18" = '17CoinExTokenCET'
    "E41" = '40KickTokenKICKBestin24h'
    "E44" = '43LocalTradersLCT'
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

Write-Host "Applied cryptos.xlsx symbol-list update"
